$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B-column values for 24.03.2025 (target rows 2-97), by quarter index 0..95
$bDay1 = @(5152,5127,5090,5049,5035,5015,4989,4956,4965,4972,4968,4985,5001,5022,5059,5103,5246,5333,5420,5528,5758,5863,6016,6211,6501,6683,6812,6886,7031,7044,7069,6997,6885,6864,6755,6647,6460,6411,6290,6219,6151,6149,6123,6106,5883,5864,5883,5950,5856,5889,5858,5946,5937,5996,5978,6025,6038,6173,6247,6368,6433,6598,6716,6805,6820,6948,6984,7102,7164,7279,7423,7576,7645,7730,7795,7786,7729,7655,7608,7537,7389,7259,7156,7008,6800,6630,6496,6405,6173,6048,5954,5877,5831,5734,5728,5644)
# New B-column values for 25.03.2025 (target rows 98-193), by quarter index 0..95
$bDay2 = @(5618,5550,5514,5488,5447,5414,5411,5415,5330,5310,5307,5318,5339,5370,5422,5464,5528,5592,5669,5744,5856,5917,6039,6176,6367,6444,6526,6537,6608,6527,6528,6475,6330,6286,6184,6101,5946,5894,5772,5729,5592,5555,5515,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

# Day 1 block: rows 2-97 -> 24.03.2025, date serial base 45740
$base1 = 45740
$prefix1 = "24.03.2025"
for ($q = 0; $q -lt 96; $q++) {
    $row = 2 + $q
    $ws.Cells.Item($row, 1).Value = [Math]::Round($base1 + $q/96, 11)
    $ws.Cells.Item($row, 2).Value = $bDay1[$q]
    $ws.Cells.Item($row, 3).Value = $q + 1
    $ws.Cells.Item($row, 4).Value = "$prefix1$($q + 1)"
}

# Day 2 block: rows 98-193 -> 25.03.2025, date serial base 45741
$base2 = 45741
$prefix2 = "25.03.2025"
for ($q = 0; $q -lt 96; $q++) {
    $row = 98 + $q
    $ws.Cells.Item($row, 1).Value = [Math]::Round($base2 + $q/96, 11)
    $ws.Cells.Item($row, 2).Value = $bDay2[$q]
    $ws.Cells.Item($row, 3).Value = $q + 1
    $ws.Cells.Item($row, 4).Value = "$prefix2$($q + 1)"
}
